$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.512.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.825.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.08%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'316.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5167"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.77%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3869"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.08412"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +9.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.120"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.84%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'41.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.85%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.825.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'94.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001132"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.74%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06641"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.55%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.51%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'6.078"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.29%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'28.563.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.02%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'21.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.44%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'159.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.60%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.034.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.08%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.402"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.11%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'125.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.02%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.1096"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.56%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.48%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.07751"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +9.52%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.745"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.676"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Algorand"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.2228"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.38%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'VeChain"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.02380"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.277"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.95%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.747"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.19%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.6439"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'11.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.192"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'13.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.84%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.6194"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.33%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.797"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.18%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'127.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.73%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.23%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'74.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.12%  "
$ws.Range("E51").Style = "Normal"
